$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("C2").Value = -1.478268366078623
$ws.Range("D2").Value = 0.1535113371196473

$ws.Range("C3").Value = 0.009527905428311093
$ws.Range("D3").Value = 0.9924838202041442

$ws.Range("C4").Value = -0.8109820795995006
$ws.Range("D4").Value = 0.4260629262353477

$ws.Range("C5").Value = -1.297667969103662
$ws.Range("D5").Value = 0.2078436233020673

$ws.Range("C6").Value = 1.277244731243567
$ws.Range("D6").Value = 0.2148315205972382

$ws.Range("C7").Value = 0.8389077716925786
$ws.Range("D7").Value = 0.4105442349659554

$ws.Range("C8").Value = 0.4509866006628294
$ws.Range("D8").Value = 0.6564129599905915

$ws.Range("C9").Value = -0.9237053280916568
$ws.Range("D9").Value = 0.3656607189909526

$ws.Range("C10").Value = -1.075413707262149
$ws.Range("D10").Value = 0.2938426642926275

$ws.Range("C11").Value = -0.4536356766284825
$ws.Range("D11").Value = 0.6545342010838859
